$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")
$checklist = $wb.Worksheets.Item("Checklist")

# --- Row 19: height changes from 45 to 53.25 (and becomes a custom height) ---
$ws.Rows.Item(19).RowHeight = 53.25

# --- Row 20: "Legend for Category field" ---
$ws.Rows.Item(20).RowHeight = 83.25
$ws.Range("A20:E20").WrapText = $true
$ws.Range("A20").Value = 15
$ws.Range("B20").Value = "Legend for Category field"
$ws.Range("C20").Value = "Display legend for one of the Category field"
$ws.Range("D20").Value = @"
1. Drag 'ProjectNane' in 'Legend' field.
2.  Go to 'Formatting pane'
3. Go to 'Bar Formatting'
4. Update 'Color' for all bars
 
"@

# --- Row 21: "Enable Hierarchy layout" ---
$ws.Rows.Item(21).RowHeight = 89.25
$ws.Range("A21:E21").WrapText = $true
$ws.Range("A21").Value = 16
$ws.Range("B21").Value = "Enable Hierarchy layout"
$ws.Range("C21").Value = "Hierarchy is enabled for Category"
$ws.Range("D21").Value = @"
1. Go to 'Formatting pane'
2.  Go to 'Column labels'
3. Enable toggle for 'Hierarchial layout'
 
"@
$ws.Range("E21").Value = "1.Hierarchial layout is enabled for Categories"

# E20 is filled in after row 21 to match the original authoring order
# (and therefore the shared-string table order).
$ws.Range("E20").Value = @"
1. Legends for 'ProjectName' will appear
4. 'Color' will be updated for the legends
"@

# --- Row 22: "Bookmarks" ---
$ws.Rows.Item(22).RowHeight = 95.25
$ws.Range("C22:E22").WrapText = $true
$ws.Range("A22").Value = 15
$ws.Range("B22").Value = "Bookmarks"
$ws.Range("C22").Value = @"
Check whether bookmarks feature works
(Note: Won’t work when Hierarchy layout is ON)
"@
$ws.Range("D22").Value = @"
1. Go to View and turn on Bookmarks Pane
2. In the visual, perform selections
3. In the boomarks pane, add a new bookmark such that selections are retained
4. Now change selections & click on the saved bookmark
"@
$ws.Range("E22").Value = @"
1. Bookmarks Pane will be visible on the left
2. Visual will update according to selections
3. In the boomarks pane, a new entry of the bookmark will come
4. The selection state saved in bookmark will be restored in the visual
"@
$ws.Range("F22").Value = "Pass"
$ws.Range("G22").Value = "Pass"
$ws.Range("H22").Value = "Pass"
$ws.Range("I22").Value = "Pass"

# --- Cursor/selection moves to E29 on the BVT sheet ---
$ws.Range("E29").Select()

# --- Checklist sheet: add portrait page setup ---
$checklist.PageSetup.Orientation = 1
